$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1730.5
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 2468.8
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 2468.8
$ws.Range("M2").Value = -387
$ws.Range("N2").Value = -2694.8

$ws.Range("H17").Value = 2442.2173
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2658.55
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 7975.650000000001
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -8311.650000000001

$ws.Range("H33").Value = 135.09091
$ws.Range("I33").Value = 112.28571
$ws.Range("J33").Value = 175
$ws.Range("K33").Value = 112.28571
$ws.Range("L33").Value = 175
$ws.Range("M33").Value = 116.71429
$ws.Range("N33").Value = -633

$ws.Range("H64").Value = 5135
$ws.Range("I64").Value = 3689
$ws.Range("J64").Value = 8750
$ws.Range("K64").Value = 3689
$ws.Range("L64").Value = 8750
$ws.Range("M64").Value = -3441
$ws.Range("N64").Value = -9246

$ws.Range("H67").Value = 5135
$ws.Range("I67").Value = 3689
$ws.Range("J67").Value = 8750
$ws.Range("K67").Value = 3689
$ws.Range("L67").Value = 8750
$ws.Range("M67").Value = -2831
$ws.Range("N67").Value = -10466

$ws.Range("H98").Value = 976.4
$ws.Range("I98").Value = 758.6667
$ws.Range("K98").Value = 758.6667
$ws.Range("M98").Value = 739.3333

$ws.Range("H122").Value = 976.4
$ws.Range("I122").Value = 758.6667
$ws.Range("K122").Value = 2276.0001
$ws.Range("M122").Value = 173.9998999999998

$ws.Range("H129").Value = 1974.4166
$ws.Range("I129").Value = 1156.5
$ws.Range("K129").Value = 3469.5
$ws.Range("M129").Value = 1530.5

$ws.Range("H132").Value = 20136.545
$ws.Range("I132").Value = 20136.545
$ws.Range("K132").Value = 60409.63499999999
$ws.Range("M132").Value = -57879.63499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 341.2
$ws.Range("I16").Value = 176.5
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 176.5
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 110.5
$ws.Range("N16").Value = -1574

$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 3500
$ws.Range("K46").Value = 3500
$ws.Range("M46").Value = -3181

$ws.Range("H61").Value = 3190.1
$ws.Range("I61").Value = 3190.1
$ws.Range("K61").Value = 3190.1
$ws.Range("M61").Value = -2978.1

$ws.Range("H110").Value = 1917.2222
$ws.Range("I110").Value = 1308.0769
$ws.Range("K110").Value = 1308.0769
$ws.Range("M110").Value = 736.9231

$ws.Range("H132").Value = 2577.4
$ws.Range("I132").Value = 2577.4
$ws.Range("K132").Value = 7732.200000000001
$ws.Range("M132").Value = -5202.200000000001

$ws.Range("H133").Value = 37666.668
$ws.Range("I133").Value = 23000
$ws.Range("J133").Value = 45000
$ws.Range("K133").Value = 23000
$ws.Range("L133").Value = 45000
$ws.Range("M133").Value = -20470
$ws.Range("N133").Value = -50060

$ws.Range("H136").Value = 3190.1
$ws.Range("I136").Value = 3190.1
$ws.Range("K136").Value = 9570.3
$ws.Range("M136").Value = -7020.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3638.111
$ws.Range("I86").Value = 1959.9231
$ws.Range("K86").Value = 1959.9231
$ws.Range("M86").Value = -836.9231

$ws.Range("H89").Value = 3638.111
$ws.Range("I89").Value = 1959.9231
$ws.Range("K89").Value = 9799.6155
$ws.Range("M89").Value = -4183.6155

$ws.Range("H134").Value = 830.2857
$ws.Range("I134").Value = 802
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 2406
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = 129
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5014.972
$ws.Range("I31").Value = 1905.7778
$ws.Range("J31").Value = 8124.1665
$ws.Range("K31").Value = 1905.7778
$ws.Range("L31").Value = 8124.1665
$ws.Range("M31").Value = -1610.7778
$ws.Range("N31").Value = -8714.1665

$ws.Range("H34").Value = 5014.972
$ws.Range("I34").Value = 1905.7778
$ws.Range("J34").Value = 8124.1665
$ws.Range("K34").Value = 1905.7778
$ws.Range("L34").Value = 8124.1665
$ws.Range("M34").Value = -1703.7778
$ws.Range("N34").Value = -8528.1665

$ws.Range("H59").Value = 35153
$ws.Range("J59").Value = 48921.668
$ws.Range("L59").Value = 48921.668
$ws.Range("N59").Value = -51211.668

$ws.Range("H108").Value = 48969.555
$ws.Range("I108").Value = 21312.5
$ws.Range("J108").Value = 56871.57
$ws.Range("K108").Value = 21312.5
$ws.Range("L108").Value = 56871.57
$ws.Range("M108").Value = -17472.5
$ws.Range("N108").Value = -64551.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 7000
$ws.Range("J82").Value = 7000
$ws.Range("L82").Value = 21000
$ws.Range("N82").Value = -21812

$ws.Range("H85").Value = 7000
$ws.Range("J85").Value = 7000
$ws.Range("L85").Value = 21000
$ws.Range("N85").Value = -23808

$ws.Range("H131").Value = 4177
$ws.Range("I131").Value = 3765
$ws.Range("J131").Value = 4314.3335
$ws.Range("K131").Value = 11295
$ws.Range("L131").Value = 12943.0005
$ws.Range("M131").Value = -6255
$ws.Range("N131").Value = -23023.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 52297.332
$ws.Range("I57").Value = 8666.667
$ws.Range("J57").Value = 95928
$ws.Range("K57").Value = 8666.667
$ws.Range("L57").Value = 95928
$ws.Range("M57").Value = -7846.666999999999
$ws.Range("N57").Value = -97568

$ws.Range("H132").Value = 129693.125
$ws.Range("I132").Value = 252886.25
$ws.Range("K132").Value = 758658.75
$ws.Range("M132").Value = -756128.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3716.5334
$ws.Range("I7").Value = 2750.5386
$ws.Range("K7").Value = 2750.5386
$ws.Range("M7").Value = -2638.5386

$ws.Range("H22").Value = 869.55554
$ws.Range("J22").Value = 1225
$ws.Range("L22").Value = 1225
$ws.Range("N22").Value = -1815

$ws.Range("H27").Value = 869.55554
$ws.Range("J27").Value = 1225
$ws.Range("L27").Value = 1225
$ws.Range("N27").Value = -1439

$ws.Range("H126").Value = 3716.5334
$ws.Range("I126").Value = 2750.5386
$ws.Range("K126").Value = 8251.6158
$ws.Range("M126").Value = -5781.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 49999
$ws.Range("J47").Value = 49999
$ws.Range("L47").Value = 49999
$ws.Range("N47").Value = -51143

$ws.Range("H96").Value = 1163.625
$ws.Range("I96").Value = 1103.8182
$ws.Range("J96").Value = 1295.2
$ws.Range("K96").Value = 1103.8182
$ws.Range("L96").Value = 1295.2
$ws.Range("M96").Value = 269.1818000000001
$ws.Range("N96").Value = -4041.2

$ws.Range("H100").Value = 1346.4

$ws.Range("H107").Value = 775.53845
$ws.Range("I107").Value = 497.57144
$ws.Range("K107").Value = 1492.71432
$ws.Range("M107").Value = 427.28568

$ws.Range("H132").Value = 1641.5883
$ws.Range("I132").Value = 1207.1333
$ws.Range("K132").Value = 3621.3999
$ws.Range("M132").Value = -1091.3999
